$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- data edits -------------------------------------------------------
# C11: 14 -> 15
$ws.Range("C11").Value = 15

# New row 13 content (order matters so shared-string indices line up:
# A13 -> 55, C13 -> 56, C12 -> 57)
$ws.Range("A13").Value = "cifar10_2_split_label_Conv4_default"
$ws.Range("C13").Value = "BAD"
$ws.Range("C12").Value = "Not that good"
$ws.Range("B13").Value = 85

# C12 / C13 pick up the same "right aligned number" style already used
# by the rest of column C.
$ws.Range("C12").Style = $ws.Range("C11").Style
$ws.Range("C13").Style = $ws.Range("C11").Style

# A13 gets a new custom font (JetBrains Mono, 9.8pt, light grey-blue).
$ws.Range("A13").Font.Name = "JetBrains Mono"
$ws.Range("A13").Font.Size = 9.8
$ws.Range("A13").Font.Color = 13023145

# --- row height ---------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 20

# --- view state -----------------------------------------------------
$ws.Range("C13").Select()
$excel.ActiveWindow.Zoom = 150
